$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new data row at row 9, duplicating the formatting of row 8 ---
# Copy the formatting (styles) of row 8 into row 9 first (keeps exact style ids),
# then copy the values of row 8 into row 9 as well.
$ws.Range("A8:AH8").Copy()
$ws.Range("A9:AH9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9:AH9").PasteSpecial(-4163)  # xlPasteValues

# Give row 9 the same row height as the other data rows (63.75)
$ws.Rows.Item(9).RowHeight = 63.75

# Adjust the few cells that differ from row 8: a new "Tutoria" session name
# and new start/end dates.
$ws.Range("D9").Value = "Tutoria53"
$ws.Range("G9").Value = 45619
$ws.Range("H9").Value = 45619

# --- Scroll the sheet view back to A1 (drop the stale topLeftCell="F1") ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1

# --- Clear the stray empty cells W3:AH3 (no value, no style) ---
$ws.Range("W3:AH3").Clear()
